$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2019-12-31 00:00:00"
$ws.Range("O2").Value = 52552399.44
$ws.Range("P2").Value = 212750230.69
$ws.Range("Q2").Value = 154372124.18
$ws.Range("R2").Value = -0.4959141338
$ws.Range("S2").Value = 131129079.18
$ws.Range("T2").Value = 131129079.18
$ws.Range("U2").Value = -0.7591002893
$ws.Range("V2").Value = 2598574.24
$ws.Range("W2").Value = 9923438.08
$ws.Range("X2").Value = -933442.34
$ws.Range("Y2").Value = 60648277.76
$ws.Range("Z2").Value = 60647538.2
$ws.Range("AA2").Value = 8095138.76
$ws.Range("AG2").Value = 1208307.07
$ws.Range("AP2").Value = -4.3118437351
$ws.Range("AQ2").Value = -10.719934104288
$ws.Range("AR2").Value = -14.158800802414
$ws.Range("AS2").Value = 50254322.97
$ws.Range("AT2").Value = -11.40106912436
